# Auto-generated edit script to apply diff changes to Kujata_Profits workbook
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 727.625
$ws.Range("I33").Value = 763
$ws.Range("J33").Value = 480
$ws.Range("K33").Value = 763
$ws.Range("L33").Value = 480
$ws.Range("M33").Value = -534
$ws.Range("N33").Value = -938
# Row 98
$ws.Range("H98").Value = 10359.6875
$ws.Range("I98").Value = 5039.615
$ws.Range("J98").Value = 33413.332
$ws.Range("K98").Value = 5039.615
$ws.Range("L98").Value = 33413.332
$ws.Range("M98").Value = -3541.615
$ws.Range("N98").Value = -36409.332
# Row 112
$ws.Range("H112").Value = 2556.027
$ws.Range("I112").Value = 768.4286
$ws.Range("J112").Value = 2973.1333
$ws.Range("K112").Value = 2305.2858
$ws.Range("L112").Value = 8919.3999
$ws.Range("M112").Value = -1197.2858
$ws.Range("N112").Value = -11135.3999
# Row 116
$ws.Range("H116").Value = 3691.4348
$ws.Range("J116").Value = 4100.222
$ws.Range("L116").Value = 4100.222
$ws.Range("N116").Value = -10984.222
# Row 122
$ws.Range("H122").Value = 10359.6875
$ws.Range("I122").Value = 5039.615
$ws.Range("J122").Value = 33413.332
$ws.Range("K122").Value = 15118.845
$ws.Range("L122").Value = 100239.996
$ws.Range("M122").Value = -12668.845
$ws.Range("N122").Value = -105139.996
# Row 125
$ws.Range("H125").Value = 1732.875
$ws.Range("I125").Value = 1656.75
$ws.Range("J125").Value = 1809
$ws.Range("K125").Value = 14910.75
$ws.Range("L125").Value = 16281
$ws.Range("M125").Value = -12450.75
$ws.Range("N125").Value = -21201
# Row 129
$ws.Range("H129").Value = 849.2
$ws.Range("J129").Value = 894.2826
$ws.Range("L129").Value = 2682.8478
$ws.Range("N129").Value = -12682.8478
# Row 132
$ws.Range("H132").Value = 12353732
$ws.Range("I132").Value = 16673575
$ws.Range("J132").Value = 11322.857
$ws.Range("K132").Value = 50020725
$ws.Range("L132").Value = 33968.571
$ws.Range("M132").Value = -50018195
$ws.Range("N132").Value = -39028.571
# Row 138
$ws.Range("H138").Value = 1530.8586
$ws.Range("J138").Value = 2023.8413
$ws.Range("L138").Value = 6071.5239
$ws.Range("N138").Value = -16351.5239

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3116.4866
$ws.Range("I32").Value = 2902.9033
$ws.Range("K32").Value = 2902.9033
$ws.Range("M32").Value = -2615.9033
# Row 132
$ws.Range("H132").Value = 1463.875
$ws.Range("I132").Value = 1173.0714
$ws.Range("K132").Value = 3519.2142
$ws.Range("M132").Value = -989.2142000000003

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 38462750
$ws.Range("I99").Value = 55556690
$ws.Range("J99").Value = 1400
$ws.Range("K99").Value = 55556690
$ws.Range("L99").Value = 1400
$ws.Range("M99").Value = -55555192
$ws.Range("N99").Value = -4396
# Row 134
$ws.Range("H134").Value = 3802.182
$ws.Range("J134").Value = 13169
$ws.Range("L134").Value = 39507
$ws.Range("N134").Value = -44577

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1154.7273
$ws.Range("I31").Value = 1141.6981
$ws.Range("J31").Value = 1500
$ws.Range("K31").Value = 1141.6981
$ws.Range("L31").Value = 1500
$ws.Range("M31").Value = -846.6981000000001
$ws.Range("N31").Value = -2090
# Row 34
$ws.Range("H34").Value = 1154.7273
$ws.Range("I34").Value = 1141.6981
$ws.Range("J34").Value = 1500
$ws.Range("K34").Value = 1141.6981
$ws.Range("L34").Value = 1500
$ws.Range("M34").Value = -939.6981000000001
$ws.Range("N34").Value = -1904
# Row 134
$ws.Range("H134").Value = 1545.7273
$ws.Range("I134").Value = 1244.6
$ws.Range("J134").Value = 1796.6666
$ws.Range("K134").Value = 3733.8
$ws.Range("L134").Value = 5389.9998
$ws.Range("M134").Value = -1198.8
$ws.Range("N134").Value = -10459.9998

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 32
$ws.Range("H32").Value = 2367.3333
$ws.Range("I32").Value = 2367.3333
$ws.Range("K32").Value = 7101.999899999999
$ws.Range("M32").Value = -6818.999899999999
# Row 104
$ws.Range("H104").Value = 4022.2
$ws.Range("J104").Value = 5281
$ws.Range("L104").Value = 15843
$ws.Range("N104").Value = -21085
# Row 131
$ws.Range("H131").Value = 10418656
$ws.Range("J131").Value = 2114.7642
$ws.Range("L131").Value = 6344.292600000001
$ws.Range("N131").Value = -16424.2926

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 1966.3572
$ws.Range("I122").Value = 2155.5557
$ws.Range("J122").Value = 1625.8
$ws.Range("K122").Value = 6466.6671
$ws.Range("L122").Value = 4877.4
$ws.Range("M122").Value = -4016.6671
$ws.Range("N122").Value = -9777.4
# Row 132
$ws.Range("H132").Value = 2199.1667
$ws.Range("I132").Value = 1650.7059
$ws.Range("K132").Value = 4952.1177
$ws.Range("M132").Value = -2422.1177

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 696.38464
$ws.Range("I16").Value = 696.38464
$ws.Range("K16").Value = 696.38464
$ws.Range("M16").Value = -526.38464
# Row 46
$ws.Range("H46").Value = 1937.125
$ws.Range("I46").Value = 899.6
$ws.Range("K46").Value = 899.6
$ws.Range("M46").Value = -711.6
# Row 55
$ws.Range("H55").Value = 222.1
$ws.Range("I55").Value = 187.5625
$ws.Range("J55").Value = 261.57144
$ws.Range("K55").Value = 187.5625
$ws.Range("L55").Value = 261.57144
$ws.Range("M55").Value = -14.5625
$ws.Range("N55").Value = -607.5714399999999
# Row 56
$ws.Range("H56").Value = 9200
# Row 58
$ws.Range("H58").Value = 2000
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
# Row 122
$ws.Range("H122").Value = 10121341
$ws.Range("I122").Value = 20240194
$ws.Range("J122").Value = 2488.7144
$ws.Range("K122").Value = 60720582
$ws.Range("L122").Value = 7466.1432
$ws.Range("M122").Value = -60718132
$ws.Range("N122").Value = -12366.1432
# Row 132
$ws.Range("H132").Value = 17935.816
$ws.Range("I132").Value = 814.6818
$ws.Range("K132").Value = 2444.0454
$ws.Range("M132").Value = 85.95460000000003
# Row 136
$ws.Range("H136").Value = 1576.7142
$ws.Range("I136").Value = 1612.4286
$ws.Range("J136").Value = 1541
$ws.Range("K136").Value = 4837.2858
$ws.Range("L136").Value = 4623
$ws.Range("M136").Value = -2287.2858
$ws.Range("N136").Value = -9723

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 250
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 250
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 500
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -2622
# Row 84
$ws.Range("H84").Value = 250
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 250
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 2500
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -13108
# Row 132
$ws.Range("H132").Value = 1426.1904
$ws.Range("I132").Value = 876.2143
$ws.Range("K132").Value = 2628.6429
$ws.Range("M132").Value = -98.64289999999983
